# Add a new row (53) to the DSA tracker sheet for the LeetCode 75
# "Determine If Two Strings Are Close" hashing question, and wire up its
# source-link hyperlink, mirroring the existing table rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$question = "1657. Determine If Two Strings Are Close"
$difficulty = "Medium"
$pattern = "Hashing"
$notes = "The crux is to understand the permutation and combination principle. The order does not matter, only if the set of items and their frequencies matter. You take counters of counters, or buckets. Collect the frequencies in maps, check if they have the same characters, and then check if the frequencies are equal. It is easiest to sort after the char check, then do a freq check. You can return freq1 == freq2 after sort."
$url = "https://leetcode.com/problems/determine-if-two-strings-are-close/solutions/4561223/beats-99-46-users-c-java-python-javascript-explained/?envType=study-plan-v2&envId=leetcode-75 "

$ws.Range("A53").Value = $question
$ws.Range("B53").Value = $difficulty
$ws.Range("B53").Style = $ws.Range("B52").Style
$ws.Range("C53").Value = $pattern
$ws.Range("D53").Value = $notes
$ws.Range("E53").Value = $url

$ws.Hyperlinks.Add($ws.Range("E53"), $url, "", "", $url)
$ws.Range("E53").Style = $ws.Range("E52").Style

$ws.Range("D54").Select()
